$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 77.90000000000001
$ws.Range("I15").Value = 77.90000000000001
$ws.Range("K15").Value = 233.7
$ws.Range("M15").Value = -64.70000000000002
$ws.Range("H59").Value = 1080.875
$ws.Range("J59").Value = 1080.875
$ws.Range("L59").Value = 3242.625
$ws.Range("N59").Value = -4356.625
$ws.Range("H64").Value = 4494.3335
$ws.Range("I64").Value = 4176.4
$ws.Range("J64").Value = 4891.75
$ws.Range("K64").Value = 4176.4
$ws.Range("L64").Value = 4891.75
$ws.Range("M64").Value = -3928.4
$ws.Range("N64").Value = -5387.75
$ws.Range("H67").Value = 4494.3335
$ws.Range("I67").Value = 4176.4
$ws.Range("J67").Value = 4891.75
$ws.Range("K67").Value = 4176.4
$ws.Range("L67").Value = 4891.75
$ws.Range("M67").Value = -3318.4
$ws.Range("N67").Value = -6607.75
$ws.Range("H138").Value = 4208.05
$ws.Range("I138").Value = 2240.7812
$ws.Range("J138").Value = 5133.8237
$ws.Range("K138").Value = 6722.3436
$ws.Range("L138").Value = 15401.4711
$ws.Range("M138").Value = -1582.3436
$ws.Range("N138").Value = -25681.4711
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1220461.6
$ws.Range("I32").Value = 12738.091
$ws.Range("K32").Value = 12738.091
$ws.Range("M32").Value = -12451.091
$ws.Range("H45").Value = 2665.3142
$ws.Range("I45").Value = 2391.037
$ws.Range("J45").Value = 3591
$ws.Range("K45").Value = 2391.037
$ws.Range("L45").Value = 3591
$ws.Range("M45").Value = -2014.037
$ws.Range("N45").Value = -4345
$ws.Range("H61").Value = 5553.8887
$ws.Range("I61").Value = 6236.7144
$ws.Range("J61").Value = 3164
$ws.Range("K61").Value = 6236.7144
$ws.Range("L61").Value = 3164
$ws.Range("M61").Value = -6024.7144
$ws.Range("N61").Value = -3588
$ws.Range("H63").Value = 5637.5
$ws.Range("I63").Value = 1766.6666
$ws.Range("J63").Value = 7960
$ws.Range("K63").Value = 1766.6666
$ws.Range("L63").Value = 7960
$ws.Range("M63").Value = -1080.6666
$ws.Range("N63").Value = -9332
$ws.Range("H66").Value = 5637.5
$ws.Range("I66").Value = 1766.6666
$ws.Range("J66").Value = 7960
$ws.Range("K66").Value = 8833.333000000001
$ws.Range("L66").Value = 39800
$ws.Range("M66").Value = -5401.333000000001
$ws.Range("N66").Value = -46664
$ws.Range("H74").Value = 946.1
$ws.Range("I74").Value = 939.3333
$ws.Range("K74").Value = 939.3333
$ws.Range("M74").Value = -65.33330000000001
$ws.Range("H77").Value = 946.1
$ws.Range("I77").Value = 939.3333
$ws.Range("K77").Value = 4696.6665
$ws.Range("M77").Value = -328.6665000000003
$ws.Range("H102").Value = 2606.125
$ws.Range("I102").Value = 2147.0908
$ws.Range("J102").Value = 3616
$ws.Range("K102").Value = 2147.0908
$ws.Range("L102").Value = 3616
$ws.Range("M102").Value = -525.0907999999999
$ws.Range("N102").Value = -6860
$ws.Range("H136").Value = 5553.8887
$ws.Range("I136").Value = 6236.7144
$ws.Range("J136").Value = 3164
$ws.Range("K136").Value = 18710.1432
$ws.Range("L136").Value = 9492
$ws.Range("M136").Value = -16160.1432
$ws.Range("N136").Value = -14592
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H86").Value = 4654.231
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 6313.125
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 6313.125
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -8559.125
$ws.Range("H89").Value = 4654.231
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 6313.125
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 31565.625
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -42797.625
$ws.Range("H105").Value = 1957.7097
$ws.Range("I105").Value = 1975.4828
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1975.4828
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = -228.4828
$ws.Range("N105").Value = -5194
$ws.Range("H134").Value = 12241.23
$ws.Range("I134").Value = 905.1429000000001
$ws.Range("J134").Value = 25466.666
$ws.Range("K134").Value = 2715.4287
$ws.Range("L134").Value = 76399.99800000001
$ws.Range("M134").Value = -180.4287000000004
$ws.Range("N134").Value = -81469.99800000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 594.625
$ws.Range("I16").Value = 594.625
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 594.625
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -307.625
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 11088.594
$ws.Range("I31").Value = 3771.0908
$ws.Range("J31").Value = 18878.193
$ws.Range("K31").Value = 3771.0908
$ws.Range("L31").Value = 18878.193
$ws.Range("M31").Value = -3476.0908
$ws.Range("N31").Value = -19468.193
$ws.Range("H34").Value = 11088.594
$ws.Range("I34").Value = 3771.0908
$ws.Range("J34").Value = 18878.193
$ws.Range("K34").Value = 3771.0908
$ws.Range("L34").Value = 18878.193
$ws.Range("M34").Value = -3569.0908
$ws.Range("N34").Value = -19282.193
$ws.Range("H50").Value = 8022
$ws.Range("J50").Value = 8022
$ws.Range("L50").Value = 8022
$ws.Range("N50").Value = -9272
$ws.Range("H51").Value = 8784.857
$ws.Range("J51").Value = 9415.666999999999
$ws.Range("L51").Value = 9415.666999999999
$ws.Range("N51").Value = -10887.667
$ws.Range("H59").Value = 13463
$ws.Range("J59").Value = 13463
$ws.Range("L59").Value = 13463
$ws.Range("N59").Value = -15753
$ws.Range("H61").Value = 8784.857
$ws.Range("J61").Value = 9415.666999999999
$ws.Range("L61").Value = 9415.666999999999
$ws.Range("N61").Value = -10111.667
$ws.Range("H62").Value = 4357.5
$ws.Range("I62").Value = 4357.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4357.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3733.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4357.5
$ws.Range("I65").Value = 4357.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21787.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18667.5
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 42884.42
$ws.Range("I86").Value = 8949.6
$ws.Range("J86").Value = 80589.78
$ws.Range("K86").Value = 8949.6
$ws.Range("L86").Value = 80589.78
$ws.Range("M86").Value = -7826.6
$ws.Range("N86").Value = -82835.78
$ws.Range("H89").Value = 42884.42
$ws.Range("I89").Value = 8949.6
$ws.Range("J89").Value = 80589.78
$ws.Range("K89").Value = 44748
$ws.Range("L89").Value = 402948.9
$ws.Range("M89").Value = -39132
$ws.Range("N89").Value = -414180.9
$ws.Range("H99").Value = 1153.5714
$ws.Range("I99").Value = 1118.5186
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 1118.5186
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = 379.4813999999999
$ws.Range("N99").Value = -5096
$ws.Range("H113").Value = 594.625
$ws.Range("I113").Value = 594.625
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 594.625
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1575.375
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 1153.5714
$ws.Range("I126").Value = 1118.5186
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 3355.5558
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -885.5558000000001
$ws.Range("N126").Value = -11240
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 68333.57000000001
$ws.Range("I121").Value = 6281.6
$ws.Range("J121").Value = 80743.96000000001
$ws.Range("K121").Value = 18844.8
$ws.Range("L121").Value = 242231.88
$ws.Range("M121").Value = -17534.8
$ws.Range("N121").Value = -244851.88
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9050.311
$ws.Range("I70").Value = 8298.5
$ws.Range("J70").Value = 10721
$ws.Range("K70").Value = 8298.5
$ws.Range("L70").Value = 10721
$ws.Range("M70").Value = -8028.5
$ws.Range("N70").Value = -11261
$ws.Range("H73").Value = 9050.311
$ws.Range("I73").Value = 8298.5
$ws.Range("J73").Value = 10721
$ws.Range("K73").Value = 8298.5
$ws.Range("L73").Value = 10721
$ws.Range("M73").Value = -7362.5
$ws.Range("N73").Value = -12593
$ws.Range("H80").Value = 72308.94
$ws.Range("I80").Value = 113794.3
$ws.Range("J80").Value = 3166.6667
$ws.Range("K80").Value = 113794.3
$ws.Range("L80").Value = 3166.6667
$ws.Range("M80").Value = -112796.3
$ws.Range("N80").Value = -5162.6667
$ws.Range("H83").Value = 72308.94
$ws.Range("I83").Value = 113794.3
$ws.Range("J83").Value = 3166.6667
$ws.Range("K83").Value = 568971.5
$ws.Range("L83").Value = 15833.3335
$ws.Range("M83").Value = -563979.5
$ws.Range("N83").Value = -25817.3335
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 34500
$ws.Range("J108").Value = 34500
$ws.Range("L108").Value = 34500
$ws.Range("N108").Value = -42180
$ws.Range("H136").Value = 2997.9583
$ws.Range("I136").Value = 1880.2258
$ws.Range("J136").Value = 5036.1763
$ws.Range("K136").Value = 5640.6774
$ws.Range("L136").Value = 15108.5289
$ws.Range("M136").Value = -3090.6774
$ws.Range("N136").Value = -20208.5289
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 33929.4
$ws.Range("J119").Value = 33929.4
$ws.Range("L119").Value = 33929.4
$ws.Range("N119").Value = -43605.4
$ws.Range("H136").Value = 1246.0646
$ws.Range("I136").Value = 817.8333
$ws.Range("J136").Value = 2714.2856
$ws.Range("K136").Value = 2453.4999
$ws.Range("L136").Value = 8142.8568
$ws.Range("M136").Value = 96.5001000000002
$ws.Range("N136").Value = -13242.8568
